# Atualização de bases das ligas, do dia: 04-04-2024 às 23:22
# This script re-shuffles the per-row match data (everything except the
# row's sequence id in column A and the constant Div/Div Original
# Name/Date columns C:E) among groups of rows, as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the "moving" data for a match row: B (id) and F..AC
# (HomeTeam .. PL_AhUnder). Column A (sequential index) and C,D,E
# (Div, Div Original Name, Date) remain attached to the physical row.
$dataCols = @(2) + (6..29)

function Get-RowData($row) {
    $vals = @{}
    foreach ($c in $dataCols) {
        $vals[$c] = $ws.Cells.Item($row, $c).Value()
    }
    return $vals
}

function Set-RowData($row, $vals) {
    foreach ($c in $dataCols) {
        $ws.Cells.Item($row, $c).Value = $vals[$c]
    }
}

# Each group maps: destination row -> source row (data that ends up in
# the destination row is taken from the "before" state of the source row).
$groups = @(
    @{ 84 = 85; 85 = 84 },
    @{ 109 = 110; 110 = 109 },
    @{ 113 = 115; 114 = 113; 115 = 114 },
    @{ 131 = 132; 132 = 131 },
    @{ 153 = 154; 154 = 155; 155 = 156; 156 = 153 }
)

foreach ($group in $groups) {
    $rows = $group.Keys
    $before = @{}
    foreach ($row in $rows) {
        $before[$row] = Get-RowData $row
    }
    foreach ($row in $rows) {
        $srcRow = $group[$row]
        Set-RowData $row $before[$srcRow]
    }
}
